$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list - Price (column D) and Volume(1h) (column E) values
$updates = @(
    @{ Row = 2;  D = "26.101.26";     E = "  +1.29%  " },
    @{ Row = 3;  D = "1.767.08";      E = "  +1.44%  " },
    @{ Row = 4;  D = "1.001";         E = "  +0.04%  " },
    @{ Row = 5;  D = "238.25";        E = "  +0.21%  " },
    @{ Row = 6;  D = "1.000";         E = "  +0.02%  " },
    @{ Row = 7;  D = "0.5249";        E = "  +3.92%  " },
    @{ Row = 8;  D = "0.2763";        E = "  +4.25%  " },
    @{ Row = 9;  D = "40.58";         E = "  -3.22%  " },
    @{ Row = 10; D = "0.06216";       E = "  +0.95%  " },
    @{ Row = 11; D = "1.772.84";      E = "  +1.84%  " },
    @{ Row = 12; D = "16.00";         E = "  +4.30%  " },
    @{ Row = 13; D = "0.07028";       E = "  +1.54%  " },
    @{ Row = 14; D = "0.6495";        E = "  +8.56%  " },
    @{ Row = 15; D = "4.530";         E = "  +0.87%  " },
    @{ Row = 16; D = "78.47";         E = "  +2.17%  " },
    @{ Row = 17; D = "1.000";         E = "  -0.01%  " },
    @{ Row = 18; D = "0.9997";        E = "  -0.06%  " },
    @{ Row = 19; D = "26.116.02";     E = "  +1.40%  " },
    @{ Row = 20; D = "11.77";         E = "  +1.40%  " },
    @{ Row = 21; D = "0.000006785";   E = "  -0.56%  " },
    @{ Row = 22; D = "2.001.05";      E = "  +1.98%  " },
    @{ Row = 23; D = "4.094";         E = "  +1.30%  " },
    @{ Row = 24; D = "8.442";         E = "  +4.11%  " },
    @{ Row = 25; D = "5.213";         E = "  +0.64%  " },
    @{ Row = 26; D = "137.88";        E = "  -0.11%  " },
    @{ Row = 27; D = "1.491";         E = "  -1.56%  " },
    @{ Row = 28; D = "1.863";         E = "  +3.20%  " },
    @{ Row = 29; D = "15.21";         E = "  +1.53%  " },
    @{ Row = 30; D = "102.76";        E = "  -0.64%  " },
    @{ Row = 31; D = "0.08411";       E = "  +3.92%  " },
    @{ Row = 32; D = "3.746";         E = "  -0.29%  " },
    @{ Row = 33; D = "3.468";         E = "  +0.13%  " },
    @{ Row = 34; D = "0.04465";       E = "  -1.11%  " },
    @{ Row = 35; D = "2.655" },
    @{ Row = 36; D = "1.010";         E = "  +3.18%  " },
    @{ Row = 37; D = "0.6117";        E = "  +0.74%  " },
    @{ Row = 38; D = "2.763";         E = "  +4.33%  " },
    @{ Row = 39; D = "2.007";         E = "  +5.65%  " },
    @{ Row = 40; D = "0.01588";       E = "  +2.45%  " },
    @{ Row = 41; D = "1.002";         E = "  +0.16%  " },
    @{ Row = 42; D = "103.05";        E = "  +0.11%  " },
    @{ Row = 43; D = "0.3915";        E = "  +2.88%  " },
    @{ Row = 44; D = "0.7551";        E = "  +3.24%  " },
    @{ Row = 45; D = "4.961";         E = "  -2.49%  " },
    @{ Row = 46; D = "6.520";         E = "  +10.51%  " },
    @{ Row = 47; E = "  +3.17%  " },
    @{ Row = 48; D = "0.1122";        E = "  +0.82%  " },
    @{ Row = 49; D = "30.42";         E = "  +0.88%  " },
    @{ Row = 50; D = "53.05";         E = "  +1.00%  " },
    @{ Row = 51; D = "0.3484";        E = "  +0.93%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.Value = "'" + $u.E
        $cell.Style = "Normal"
    }
}
